$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("metrics_metadata")
$wsStudy = $wb.Worksheets.Item("study")

# Remove the "scenario" column (I) values from metrics_metadata (rows 2-17) -
# the workbook no longer tags rows with a baseline scenario.
$wsMeta.Range("I2:I17").ClearContents()

# Selection / active-sheet bookkeeping that Excel records when a user
# finishes editing on "metrics_metadata" (now the active tab) after having
# started on "study".
$wsStudy.Range("E57").Select() | Out-Null
$wsMeta.Activate() | Out-Null
$wsMeta.Range("M19").Select() | Out-Null
